$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 currently holds the "men_boohoo_coat_10" listing (Single Breasted
# Wool Mix Overcoat) and row 12 holds "men_boohoo_coat_11" (Skinny Fit
# Pinstripe Double Breasted Suit). The edit swaps the Image/Title/Price/URL
# content (columns B, D, E, F) between these two rows while keeping each
# row's Id (column A) and Brand (column C) unchanged - and along the way
# corrects the overcoat's price to $66.00 (it lands on row 12 after the
# swap) to match the rest of the sheet.

# Row 11 becomes the "Skinny Fit Pinstripe Double Breasted Suit" listing.
$ws.Range("B11").Value = "https://media.boohoo.com/i/boohoo/bmm65763_navy_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit"
$ws.Range("D11").Value = "SKINNY FIT PINSTRIPE DOUBLE BREASTED SUIT"
$ws.Range("F11").Value = "https://ca.boohoo.com/skinny-fit-pinstripe-double-breasted-suit/MAN03215.html"

# Row 12 becomes the "Single Breasted Wool Mix Overcoat" listing.
$ws.Range("B12").Value = "https://media.boohoo.com/i/boohoo/mzz13957_black_xl?w=900&qlt=default&fmt.jp2.qlt=70&fmt=auto&sm=fit"
$ws.Range("D12").Value = "SINGLE BREASTED WOOL MIX OVERCOAT"
$ws.Range("F12").Value = "https://ca.boohoo.com/single-breasted-wool-mix-overcoat/MZZ13957.html"

# Prices ("$40.00" / "$66.00") need to stay plain text cells (matching the
# rest of the price column), but assigning those literal strings via
# .Value gets auto-coerced to a currency number by the COM layer. Instead,
# copy the already-text price cells that hold the exact values we need, so
# the destination cells end up with the same text shared-string type and
# no new number-format style is introduced.
$ws.Range("E21").Copy($ws.Range("E11"))
$ws.Range("E10").Copy($ws.Range("E12"))
